$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 15196.4
$ws.Range("I18").Value = 15829.444
$ws.Range("K18").Value = 15829.444
$ws.Range("M18").Value = -15545.444
$ws.Range("H32").Value = 9666.666999999999
$ws.Range("J32").Value = 9666.666999999999
$ws.Range("L32").Value = 9666.666999999999
$ws.Range("N32").Value = -10318.667
$ws.Range("H34").Value = 899.6667
$ws.Range("I34").Value = 899.6667
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 899.6667
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -696.6667
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 899.6667
$ws.Range("I36").Value = 899.6667
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 899.6667
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -184.6667
$ws.Range("N36").ClearContents()
$ws.Range("H42").Value = 43
$ws.Range("I42").Value = 8
$ws.Range("J42").Value = 66.333336
$ws.Range("K42").Value = 24
$ws.Range("L42").Value = 199.000008
$ws.Range("M42").Value = 206
$ws.Range("N42").Value = -659.000008
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H135").Value = 1366.9333
$ws.Range("I135").Value = 1239.6923
$ws.Range("K135").Value = 11157.2307
$ws.Range("M135").Value = -8622.2307
$ws.Range("H137").Value = 2083
$ws.Range("I137").Value = 1999.75
$ws.Range("J137").Value = 2249.5
$ws.Range("K137").Value = 5999.25
$ws.Range("L137").Value = 6748.5
$ws.Range("M137").Value = -3449.25
$ws.Range("N137").Value = -11848.5
$ws.Range("H141").Value = 3342.875
$ws.Range("I141").Value = 3342.875
$ws.Range("K141").Value = 10028.625
$ws.Range("M141").Value = -4848.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3500
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 3500
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -5248
$ws.Range("H77").Value = 3500
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 17500
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -26236
$ws.Range("H132").Value = 1924
$ws.Range("I132").Value = 1744.75
$ws.Range("K132").Value = 5234.25
$ws.Range("M132").Value = -2704.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 398.5
$ws.Range("I64").Value = 398.5
$ws.Range("K64").Value = 398.5
$ws.Range("M64").Value = -173.5
$ws.Range("H67").Value = 398.5
$ws.Range("I67").Value = 398.5
$ws.Range("K67").Value = 398.5
$ws.Range("M67").Value = 381.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 152.25
$ws.Range("J10").Value = 1000
$ws.Range("L10").Value = 1000
$ws.Range("N10").Value = -1278
$ws.Range("H31").Value = 6305.1665
$ws.Range("I31").Value = 7199.8335
$ws.Range("J31").Value = 5410.5
$ws.Range("K31").Value = 7199.8335
$ws.Range("L31").Value = 5410.5
$ws.Range("M31").Value = -6904.8335
$ws.Range("N31").Value = -6000.5
$ws.Range("H34").Value = 6305.1665
$ws.Range("I34").Value = 7199.8335
$ws.Range("J34").Value = 5410.5
$ws.Range("K34").Value = 7199.8335
$ws.Range("L34").Value = 5410.5
$ws.Range("M34").Value = -6997.8335
$ws.Range("N34").Value = -5814.5
$ws.Range("H92").Value = 20000
$ws.Range("J92").Value = 20000
$ws.Range("L92").Value = 20000
$ws.Range("N92").Value = -24992
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H122").Value = 2112
$ws.Range("I122").Value = 2112
$ws.Range("K122").Value = 6336
$ws.Range("M122").Value = -3886
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 1071.909
$ws.Range("I132").Value = 786.625
$ws.Range("K132").Value = 2359.875
$ws.Range("M132").Value = 170.125
$ws.Range("H134").Value = 6750
$ws.Range("I134").Value = 6000
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 18000
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -15465
$ws.Range("N134").Value = -26070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 671.3333
$ws.Range("I12").Value = 7
$ws.Range("K12").Value = 21
$ws.Range("M12").Value = 152
$ws.Range("H23").Value = 767.13336
$ws.Range("I23").Value = 711.7778
$ws.Range("K23").Value = 2135.3334
$ws.Range("M23").Value = -1900.3334
$ws.Range("H107").Value = 299.5
$ws.Range("J107").Value = 500
$ws.Range("L107").Value = 1500
$ws.Range("N107").Value = -5340
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H129").Value = 2018.4
$ws.Range("I129").Value = 1764.75
$ws.Range("K129").Value = 5294.25
$ws.Range("M129").Value = -294.25
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H136").Value = 29500
$ws.Range("J136").Value = 29500
$ws.Range("L136").Value = 88500
$ws.Range("N136").Value = -93600

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1999.5
$ws.Range("I22").Value = 1999
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 1999
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -1704
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 1999.5
$ws.Range("I27").Value = 1999
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 1999
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -1892
$ws.Range("N27").Value = -2214
$ws.Range("H40").Value = 7312.5
$ws.Range("I40").Value = 7312.5
$ws.Range("K40").Value = 7312.5
$ws.Range("M40").Value = -7176.5
$ws.Range("H55").Value = 5944.3335
$ws.Range("J55").Value = 5944.3335
$ws.Range("L55").Value = 5944.3335
$ws.Range("N55").Value = -6290.3335
$ws.Range("H68").Value = 2900
$ws.Range("I68").Value = 2900
$ws.Range("K68").Value = 2900
$ws.Range("M68").Value = -2151
$ws.Range("H71").Value = 2900
$ws.Range("I71").Value = 2900
$ws.Range("K71").Value = 14500
$ws.Range("M71").Value = -10756
$ws.Range("H76").Value = 20461.834
$ws.Range("J76").Value = 20697
$ws.Range("L76").Value = 20697
$ws.Range("N76").Value = -21373
$ws.Range("H79").Value = 20461.834
$ws.Range("J79").Value = 20697
$ws.Range("L79").Value = 20697
$ws.Range("N79").Value = -23037
$ws.Range("H100").Value = 3339.9092
$ws.Range("I100").Value = 2781.6667
$ws.Range("K100").Value = 2781.6667
$ws.Range("M100").Value = -2240.6667
$ws.Range("H122").Value = 7333.3335
$ws.Range("I122").Value = 8000
$ws.Range("K122").Value = 24000
$ws.Range("M122").Value = -21550

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 40604.168
$ws.Range("J45").Value = 42075
$ws.Range("L45").Value = 42075
$ws.Range("N45").Value = -43057
$ws.Range("H81").Value = 725.3333
$ws.Range("I81").Value = 588
$ws.Range("K81").Value = 1176
$ws.Range("M81").Value = -115
$ws.Range("H84").Value = 725.3333
$ws.Range("I84").Value = 588
$ws.Range("K84").Value = 5880
$ws.Range("M84").Value = -576
$ws.Range("H100").Value = 3067.0715
$ws.Range("I100").Value = 2779.3635
$ws.Range("J100").Value = 4122
$ws.Range("K100").Value = 5558.727
$ws.Range("L100").Value = 8244
$ws.Range("M100").Value = -5017.727
$ws.Range("N100").Value = -9326
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H113").Value = 444.7
$ws.Range("I113").Value = 427.55554
$ws.Range("K113").Value = 1282.66662
$ws.Range("M113").Value = 887.33338
$ws.Range("H132").Value = 3353.8823
$ws.Range("I132").Value = 1593.25
$ws.Range("K132").Value = 4779.75
$ws.Range("M132").Value = -2249.75
